$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ----------------------------------------------------------
# B1/C1 already carry the bold/bordered header style (style index 1).
# Copy their formatting down into A1/B1 (the new header positions) before
# overwriting the text, so the header keeps its look; then drop the
# now-unused C1 cell.
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("C1").Clear()

$ws.Range("A1").Value = "Average"
$ws.Range("B1").Value = "Standard Deviation"

# --- Data rows -------------------------------------------------------------
# The old A2:A9 label cells carried the bold/bordered style too; clear the
# whole old data block first (values + formatting) so the new plain numeric
# cells don't inherit it.
$ws.Range("A2:C9").Clear()

# Replace the old 8-row label/avg/stdev table with the new 3-row
# average/standard-deviation summary.
$ws.Range("A2").Value = 53.33333333333334
$ws.Range("B2").Value = 0.4714045207910317

$ws.Range("A3").Value = 53.33333333333334
$ws.Range("B3").Value = 0.4714045207910317

$ws.Range("A4").Value = 52
$ws.Range("B4").Value = 0
